$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 180.5
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""
$ws.Range("H53").Value = 13889783
$ws.Range("J53").Value = 1194.5385
$ws.Range("L53").Value = 1194.5385
$ws.Range("N53").Value = -2468.5385
$ws.Range("H135").Value = 1744.1
$ws.Range("I135").Value = 1426.375
$ws.Range("J135").Value = 3015
$ws.Range("K135").Value = 12837.375
$ws.Range("L135").Value = 27135
$ws.Range("M135").Value = -10302.375
$ws.Range("N135").Value = -32205
$ws.Range("H137").Value = 3216.2888
$ws.Range("I137").Value = 2573.7576
$ws.Range("J137").Value = 4983.25
$ws.Range("K137").Value = 7721.2728
$ws.Range("L137").Value = 14949.75
$ws.Range("M137").Value = -5171.2728
$ws.Range("N137").Value = -20049.75
$ws.Range("H138").Value = 6061.405
$ws.Range("I138").Value = 2659.4583
$ws.Range("J138").Value = 7422.183
$ws.Range("K138").Value = 7978.374899999999
$ws.Range("L138").Value = 22266.549
$ws.Range("M138").Value = -2838.374899999999
$ws.Range("N138").Value = -32546.549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2442.06
$ws.Range("I32").Value = 2365.7173
$ws.Range("K32").Value = 2365.7173
$ws.Range("M32").Value = -2078.7173
$ws.Range("H61").Value = 3143.9092
$ws.Range("I61").Value = 2839.5789
$ws.Range("K61").Value = 2839.5789
$ws.Range("M61").Value = -2627.5789
$ws.Range("H74").Value = 2093.0217
$ws.Range("I74").Value = 2101.5715
$ws.Range("J74").Value = 2003.25
$ws.Range("K74").Value = 2101.5715
$ws.Range("L74").Value = 2003.25
$ws.Range("M74").Value = -1227.5715
$ws.Range("N74").Value = -3751.25
$ws.Range("H76").Value = 100000
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = ""
$ws.Range("H77").Value = 2093.0217
$ws.Range("I77").Value = 2101.5715
$ws.Range("J77").Value = 2003.25
$ws.Range("K77").Value = 10507.8575
$ws.Range("L77").Value = 10016.25
$ws.Range("M77").Value = -6139.8575
$ws.Range("N77").Value = -18752.25
$ws.Range("H79").Value = 100000
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = ""
$ws.Range("H98").Value = 53451.332
$ws.Range("J98").Value = 53451.332
$ws.Range("L98").Value = 53451.332
$ws.Range("N98").Value = -59441.332
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""
$ws.Range("H136").Value = 3143.9092
$ws.Range("I136").Value = 2839.5789
$ws.Range("K136").Value = 8518.736699999999
$ws.Range("M136").Value = -5968.736699999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 48056.184
$ws.Range("I31").Value = 1794.3636
$ws.Range("J31").Value = 94318
$ws.Range("K31").Value = 1794.3636
$ws.Range("L31").Value = 94318
$ws.Range("M31").Value = -1499.3636
$ws.Range("N31").Value = -94908
$ws.Range("H34").Value = 48056.184
$ws.Range("I34").Value = 1794.3636
$ws.Range("J34").Value = 94318
$ws.Range("K34").Value = 1794.3636
$ws.Range("L34").Value = 94318
$ws.Range("M34").Value = -1592.3636
$ws.Range("N34").Value = -94722
$ws.Range("H58").Value = 2666.3076
$ws.Range("J58").Value = 2728.125
$ws.Range("L58").Value = 2728.125
$ws.Range("N58").Value = -3134.125
$ws.Range("H86").Value = 11111
$ws.Range("J86").Value = 11111
$ws.Range("L86").Value = 11111
$ws.Range("N86").Value = -13357
$ws.Range("H89").Value = 11111
$ws.Range("J89").Value = 11111
$ws.Range("L89").Value = 55555
$ws.Range("N89").Value = -66787
$ws.Range("H103").Value = 3762
$ws.Range("I103").Value = 3762
$ws.Range("K103").Value = 3762
$ws.Range("M103").Value = -2590
$ws.Range("H132").Value = 4764.3335
$ws.Range("I132").Value = 2610
$ws.Range("K132").Value = 7830
$ws.Range("M132").Value = -5300
$ws.Range("H136").Value = 2666.3076
$ws.Range("J136").Value = 2728.125
$ws.Range("L136").Value = 8184.375
$ws.Range("N136").Value = -13284.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 900583.25
$ws.Range("I5").Value = 89510.44500000001
$ws.Range("J5").Value = 3333801.8
$ws.Range("K5").Value = 268531.335
$ws.Range("L5").Value = 10001405.4
$ws.Range("M5").Value = -268419.335
$ws.Range("N5").Value = -10001629.4
$ws.Range("H109").Value = 38850.38
$ws.Range("I109").Value = 2805.1667
$ws.Range("J109").Value = 48253.477
$ws.Range("K109").Value = 8415.500100000001
$ws.Range("L109").Value = 144760.431
$ws.Range("M109").Value = -7375.500100000001
$ws.Range("N109").Value = -146840.431
$ws.Range("H131").Value = 24694.113
$ws.Range("I131").Value = 78345.16
$ws.Range("J131").Value = 15516.961
$ws.Range("K131").Value = 235035.48
$ws.Range("L131").Value = 46550.883
$ws.Range("M131").Value = -229995.48
$ws.Range("N131").Value = -56630.883
$ws.Range("H135").Value = 900583.25
$ws.Range("I135").Value = 89510.44500000001
$ws.Range("J135").Value = 3333801.8
$ws.Range("K135").Value = 805594.0050000001
$ws.Range("L135").Value = 30004216.2
$ws.Range("M135").Value = -803059.0050000001
$ws.Range("N135").Value = -30009286.2
$ws.Range("H139").Value = 6254.919
$ws.Range("I139").Value = 3271
$ws.Range("J139").Value = 8528.380999999999
$ws.Range("K139").Value = 9813
$ws.Range("L139").Value = 25585.143
$ws.Range("M139").Value = -4673
$ws.Range("N139").Value = -35865.143
$ws.Range("H140").Value = 3108.9285
$ws.Range("I140").Value = 2509.1765
$ws.Range("J140").Value = 4035.818
$ws.Range("K140").Value = 7527.529500000001
$ws.Range("L140").Value = 12107.454
$ws.Range("M140").Value = -2347.529500000001
$ws.Range("N140").Value = -22467.454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 56003
$ws.Range("I19").Value = 15000
$ws.Range("K19").Value = 15000
$ws.Range("M19").Value = -14712
$ws.Range("H47").Value = 40031
$ws.Range("J47").Value = 40031
$ws.Range("L47").Value = 40031
$ws.Range("N47").Value = -41167
$ws.Range("H58").Value = 30045.889
$ws.Range("I58").Value = 30045
$ws.Range("K58").Value = 30045
$ws.Range("M58").Value = -29768
$ws.Range("H80").Value = 561168.7
$ws.Range("I80").Value = 481268.62
$ws.Range("J80").Value = 673028.8
$ws.Range("K80").Value = 481268.62
$ws.Range("L80").Value = 673028.8
$ws.Range("M80").Value = -480270.62
$ws.Range("N80").Value = -675024.8
$ws.Range("H83").Value = 561168.7
$ws.Range("I83").Value = 481268.62
$ws.Range("J83").Value = 673028.8
$ws.Range("K83").Value = 2406343.1
$ws.Range("L83").Value = 3365144
$ws.Range("M83").Value = -2401351.1
$ws.Range("N83").Value = -3375128
$ws.Range("H132").Value = 74127.664
$ws.Range("I132").Value = 8372.25
$ws.Range("K132").Value = 25116.75
$ws.Range("M132").Value = -22586.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5732.85
$ws.Range("I132").Value = 4691.615
$ws.Range("J132").Value = 7666.5713
$ws.Range("K132").Value = 14074.845
$ws.Range("L132").Value = 22999.7139
$ws.Range("M132").Value = -11544.845
$ws.Range("N132").Value = -28059.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = ""
